# "Generate Report for Handback" - refresh the localization-status report
# after a handback event: the zh-cn / de-de targets are now in sync with
# en-US, so the previous "stale handback" error clears, the handback
# timestamps move forward, and the status text updates everywhere it is
# shown (Overview sheet + per-locale detail sheets).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn detail sheet -------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("K2").Value = "2016-10-18 12:30:25"
$wsZh.Range("P2").Value = ""
$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---- de-de detail sheet -------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("K2").Value = "2016-10-18 12:30:42"
$wsDe.Range("P2").Value = ""
$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(16).ColumnWidth = 12.833333333333334
